# Tue Mar 21 13:37:29 UTC 2023 cryptos list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values must stay plain text (inline-string cells in the
# source data), so force text format, assign the literal, then drop back to
# the default 'Normal' style so no stray number-format style sticks around.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.143.22"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.12%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.818.90"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.06%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.71%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "338.14"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9982"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4255"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +11.07%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3539"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.99%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "45.84"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -2.51%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.167"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.95%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07520"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +1.50%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "23.07"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.39%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.0000"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.21%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.324"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.329"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.812.90"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.33%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.00001095"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.55%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.06697"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "82.84"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.9989"
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "17.45"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.01%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.394"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.61%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "28.156.09"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.92"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.407"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.59%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.490"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.15%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "20.75"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "155.95"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.95%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.018.16"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.331"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -6.62%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "133.77"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.46%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.080"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.59%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "6.057"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.04%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.09142"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.95%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "12.46"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.46%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.06361"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.58%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.02357"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.28%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.6690"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.52%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.284"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.2168"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.515"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.83%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.223"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.30%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.211"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.15%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "14.26"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.05%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.9989"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.6205"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.43%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.880"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "129.19"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.94%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.068"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.82%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.190"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.07127"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -5.10%  "

